$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:C1): add a thin box border, keep existing yellow fill ---
$ws.Range("A1:C1").Borders.LineStyle = 1

# --- A2 ("No" counter cell): add border, keep existing center/top alignment ---
$ws.Range("A2").Borders.LineStyle = 1

# --- B2 (question text cell): add border, keep existing left/top/wrap alignment ---
$ws.Range("B2").Borders.LineStyle = 1

# --- C2: new "File Name" value cell, vertical-top aligned, bordered ---
$ws.Range("C2").Value = "MergeInterval.java"
$ws.Range("C2").Borders.LineStyle = 1
$ws.Range("C2").VerticalAlignment = -4160

# --- Column C is now wider to fit the file name, no longer auto best-fit ---
$ws.Columns("C").ColumnWidth = 18.5

# --- Selection moves to B2 ---
$ws.Range("B2").Select() | Out-Null
